$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 986
$ws.Range("I29").Value = 986
$ws.Range("K29").Value = 2958
$ws.Range("M29").Value = -2677
$ws.Range("H69").Value = 10917.214
$ws.Range("I69").Value = 6316.3335
$ws.Range("J69").Value = 14367.875
$ws.Range("K69").Value = 18949.0005
$ws.Range("L69").Value = 43103.625
$ws.Range("M69").Value = -18075.0005
$ws.Range("N69").Value = -44851.625
$ws.Range("H72").Value = 10917.214
$ws.Range("I72").Value = 6316.3335
$ws.Range("J72").Value = 14367.875
$ws.Range("K72").Value = 56847.0015
$ws.Range("L72").Value = 129310.875
$ws.Range("M72").Value = -52479.0015
$ws.Range("N72").Value = -138046.875
$ws.Range("H92").Value = 841.7778
$ws.Range("I92").Value = 776.95
$ws.Range("J92").Value = 1027
$ws.Range("K92").Value = 776.95
$ws.Range("L92").Value = 1027
$ws.Range("M92").Value = 471.05
$ws.Range("N92").Value = -3523
$ws.Range("H93").Value = 42800.5
$ws.Range("J93").Value = 42800.5
$ws.Range("L93").Value = 42800.5
$ws.Range("N93").Value = -47792.5
$ws.Range("H95").Value = 50624
$ws.Range("J95").Value = 50624
$ws.Range("L95").Value = 50624
$ws.Range("N95").Value = -56116
$ws.Range("H112").Value = 4987.909
$ws.Range("J112").Value = 3709
$ws.Range("L112").Value = 11127
$ws.Range("N112").Value = -13343
$ws.Range("H135").Value = 1961.7391
$ws.Range("I135").Value = 1131.8823
$ws.Range("J135").Value = 4313
$ws.Range("K135").Value = 10186.9407
$ws.Range("L135").Value = 38817
$ws.Range("M135").Value = -7651.940699999999
$ws.Range("N135").Value = -43887
$ws.Range("H137").Value = 5484.4443
$ws.Range("I137").Value = 4226.067
$ws.Range("J137").Value = 7057.4165
$ws.Range("K137").Value = 12678.201
$ws.Range("L137").Value = 21172.2495
$ws.Range("M137").Value = -10128.201
$ws.Range("N137").Value = -26272.2495
$ws.Range("H138").Value = 8104.72
$ws.Range("I138").Value = 5648.5
$ws.Range("J138").Value = 8318.305
$ws.Range("K138").Value = 16945.5
$ws.Range("L138").Value = 24954.915
$ws.Range("M138").Value = -11805.5
$ws.Range("N138").Value = -35234.915
$ws.Range("H140").Value = 362500
$ws.Range("J140").Value = 362500
$ws.Range("L140").Value = 362500
$ws.Range("N140").Value = -372860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4978.324
$ws.Range("I2").Value = 3358.2188
$ws.Range("J2").Value = 15347
$ws.Range("K2").Value = 3358.2188
$ws.Range("L2").Value = 15347
$ws.Range("M2").Value = -3245.2188
$ws.Range("N2").Value = -15573
$ws.Range("H32").Value = 2030
$ws.Range("I32").Value = 533.02563
$ws.Range("J32").Value = 9327.75
$ws.Range("K32").Value = 533.02563
$ws.Range("L32").Value = 9327.75
$ws.Range("M32").Value = -246.02563
$ws.Range("N32").Value = -9901.75
$ws.Range("H45").Value = 2886.1738
$ws.Range("I45").Value = 3037.9524
$ws.Range("J45").Value = 1292.5
$ws.Range("K45").Value = 3037.9524
$ws.Range("L45").Value = 1292.5
$ws.Range("M45").Value = -2660.9524
$ws.Range("N45").Value = -2046.5
$ws.Range("H116").Value = 4978.324
$ws.Range("I116").Value = 3358.2188
$ws.Range("J116").Value = 15347
$ws.Range("K116").Value = 3358.2188
$ws.Range("L116").Value = 15347
$ws.Range("M116").Value = -1064.2188
$ws.Range("N116").Value = -19935
$ws.Range("H132").Value = 45820.957
$ws.Range("I132").Value = 7358.294
$ws.Range("J132").Value = 139230.28
$ws.Range("K132").Value = 22074.882
$ws.Range("L132").Value = 417690.84
$ws.Range("M132").Value = -19544.882
$ws.Range("N132").Value = -422750.84

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4982.3516
$ws.Range("I3").Value = 3286.7273
$ws.Range("K3").Value = 3286.7273
$ws.Range("M3").Value = -3172.7273
$ws.Range("H94").Value = 3822.6667
$ws.Range("I94").Value = 684.2941
$ws.Range("J94").Value = 9157.9
$ws.Range("K94").Value = 684.2941
$ws.Range("L94").Value = 9157.9
$ws.Range("M94").Value = -233.2941
$ws.Range("N94").Value = -10059.9
$ws.Range("H99").Value = 9230.558999999999
$ws.Range("I99").Value = 8994.767
$ws.Range("J99").Value = 10999
$ws.Range("K99").Value = 8994.767
$ws.Range("L99").Value = 10999
$ws.Range("M99").Value = -7496.767
$ws.Range("N99").Value = -13995
$ws.Range("H105").Value = 3163.4
$ws.Range("I105").Value = 2422.7307
$ws.Range("J105").Value = 7977.75
$ws.Range("K105").Value = 2422.7307
$ws.Range("L105").Value = 7977.75
$ws.Range("M105").Value = -675.7307000000001
$ws.Range("N105").Value = -11471.75
$ws.Range("H130").Value = 150000
$ws.Range("J130").Value = 150000
$ws.Range("L130").Value = 150000
$ws.Range("N130").Value = -160040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4059
$ws.Range("I31").Value = 1013.4375
$ws.Range("J31").Value = 4808.677
$ws.Range("K31").Value = 1013.4375
$ws.Range("L31").Value = 4808.677
$ws.Range("M31").Value = -718.4375
$ws.Range("N31").Value = -5398.677
$ws.Range("H34").Value = 4059
$ws.Range("I34").Value = 1013.4375
$ws.Range("J34").Value = 4808.677
$ws.Range("K34").Value = 1013.4375
$ws.Range("L34").Value = 4808.677
$ws.Range("M34").Value = -811.4375
$ws.Range("N34").Value = -5212.677
$ws.Range("H36").Value = 9000
$ws.Range("I36").Value = 9000
$ws.Range("K36").Value = 9000
$ws.Range("M36").Value = -8612
$ws.Range("H40").Value = 9000
$ws.Range("I40").Value = 9000
$ws.Range("K40").Value = 9000
$ws.Range("M40").Value = -8840
$ws.Range("H58").Value = 2983.6177
$ws.Range("I58").Value = 2271.32
$ws.Range("J58").Value = 4962.222
$ws.Range("K58").Value = 2271.32
$ws.Range("L58").Value = 4962.222
$ws.Range("M58").Value = -2068.32
$ws.Range("N58").Value = -5368.222
$ws.Range("H122").Value = 2550.0625
$ws.Range("I122").Value = 1901.375
$ws.Range("J122").Value = 3198.75
$ws.Range("K122").Value = 5704.125
$ws.Range("L122").Value = 9596.25
$ws.Range("M122").Value = -3254.125
$ws.Range("N122").Value = -14496.25
$ws.Range("H124").Value = 38663
$ws.Range("J124").Value = 38663
$ws.Range("L124").Value = 38663
$ws.Range("N124").Value = -43573
$ws.Range("H136").Value = 2983.6177
$ws.Range("I136").Value = 2271.32
$ws.Range("J136").Value = 4962.222
$ws.Range("K136").Value = 6813.960000000001
$ws.Range("L136").Value = 14886.666
$ws.Range("M136").Value = -4263.960000000001
$ws.Range("N136").Value = -19986.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 4094.4285
$ws.Range("I46").Value = 2124.75
$ws.Range("J46").Value = 4882.3
$ws.Range("K46").Value = 6374.25
$ws.Range("L46").Value = 14646.9
$ws.Range("M46").Value = -6283.25
$ws.Range("N46").Value = -14828.9
$ws.Range("H68").Value = 242509.33
$ws.Range("I68").Value = 1999.1428
$ws.Range("J68").Value = 362764.44
$ws.Range("K68").Value = 5997.428400000001
$ws.Range("L68").Value = 1088293.32
$ws.Range("M68").Value = -5186.428400000001
$ws.Range("N68").Value = -1089915.32
$ws.Range("H71").Value = 242509.33
$ws.Range("I71").Value = 1999.1428
$ws.Range("J71").Value = 362764.44
$ws.Range("K71").Value = 17992.2852
$ws.Range("L71").Value = 3264879.96
$ws.Range("M71").Value = -13936.2852
$ws.Range("N71").Value = -3272991.96
$ws.Range("H121").Value = 26232.5
$ws.Range("I121").Value = 1500
$ws.Range("J121").Value = 28480.908
$ws.Range("K121").Value = 4500
$ws.Range("L121").Value = 85442.724
$ws.Range("M121").Value = -3190
$ws.Range("N121").Value = -88062.724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 16667006
$ws.Range("J35").Value = 508.5
$ws.Range("L35").Value = 508.5
$ws.Range("N35").Value = -1104.5
$ws.Range("H113").Value = 9452.929
$ws.Range("I113").Value = 4668.8
$ws.Range("K113").Value = 4668.8
$ws.Range("M113").Value = -2498.8
$ws.Range("H132").Value = 4377.524
$ws.Range("I132").Value = 4267.5
$ws.Range("J132").Value = 4729.6
$ws.Range("K132").Value = 12802.5
$ws.Range("L132").Value = 14188.8
$ws.Range("M132").Value = -10272.5
$ws.Range("N132").Value = -19248.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 14116.029
$ws.Range("I61").Value = 12218.741
$ws.Range("J61").Value = 21434.143
$ws.Range("K61").Value = 12218.741
$ws.Range("L61").Value = 21434.143
$ws.Range("M61").Value = -12016.741
$ws.Range("N61").Value = -21838.143
$ws.Range("H113").Value = 14116.029
$ws.Range("I113").Value = 12218.741
$ws.Range("J113").Value = 21434.143
$ws.Range("K113").Value = 12218.741
$ws.Range("L113").Value = 21434.143
$ws.Range("M113").Value = -10048.741
$ws.Range("N113").Value = -25774.143
$ws.Range("H132").Value = 2675.8708
$ws.Range("I132").Value = 2663.5386
$ws.Range("K132").Value = 7990.6158
$ws.Range("M132").Value = -5460.6158
$ws.Range("H136").Value = 31252464
$ws.Range("I136").Value = 1707.3334
$ws.Range("J136").Value = 90913000
$ws.Range("K136").Value = 5122.0002
$ws.Range("L136").Value = 272739000
$ws.Range("M136").Value = -2572.0002
$ws.Range("N136").Value = -272744100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 7153487.5
$ws.Range("I32").Value = 7153487.5
$ws.Range("K32").Value = 7153487.5
$ws.Range("M32").Value = -7153170.5
$ws.Range("H45").Value = 16299.667
$ws.Range("J45").Value = 11999.5
$ws.Range("L45").Value = 11999.5
$ws.Range("N45").Value = -12981.5
$ws.Range("H100").Value = 1323.1666
$ws.Range("I100").Value = 1307.091
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 2614.182
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2073.182
$ws.Range("N100").Value = -4082
$ws.Range("H132").Value = 4867.5835
$ws.Range("I132").Value = 1007.2143
$ws.Range("J132").Value = 10272.1
$ws.Range("K132").Value = 3021.6429
$ws.Range("L132").Value = 30816.3
$ws.Range("M132").Value = -491.6428999999998
$ws.Range("N132").Value = -35876.3
$ws.Range("H136").Value = 2081.077
$ws.Range("I136").Value = 1462.1842
$ws.Range("J136").Value = 3760.9285
$ws.Range("K136").Value = 4386.5526
$ws.Range("L136").Value = 11282.7855
$ws.Range("M136").Value = -1836.5526
$ws.Range("N136").Value = -16382.7855
